# C5-PowerPoint.pptx edit:
#  1. Slide 6's table switches from table style {9520165D-9965-4FD6-B1FF-EAF4D9BA06F0}
#     to {D4CC77DC-F8EC-4B10-83C1-847AFAF6E699}.
#  2. The deck's theme colour palette ("ppt/theme/theme2.xml", the theme used by
#     the slide master / all slides) is swapped from the "Integral" palette to the
#     stock "Office Theme" palette (font scheme / format scheme were already
#     identical between the two themes, so only the 12 theme colours move).

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$tableSlide = $p.Slides.Item(6)
foreach ($shp in $tableSlide.Shapes) {
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{D4CC77DC-F8EC-4B10-83C1-847AFAF6E699}")
    }
}

# --- 2. Theme colours -------------------------------------------------------
# Office Theme colour scheme, in ThemeColorScheme / RGB() order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeThemeColors = @(
    0,           # dk1      000000
    16777215,    # lt1      FFFFFF
    6968388,     # dk2      44546A
    15132391,    # lt2      E7E6E6
    13998939,    # accent1  5B9BD5
    3243501,     # accent2  ED7D31
    10855845,    # accent3  A5A5A5
    49407,       # accent4  FFC000
    12874308,    # accent5  4472C4
    4697456,     # accent6  70AD47
    12673797,    # hlink    0563C1
    7491477      # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = $officeThemeColors[$i - 1]
}
